# Update transition-matrix probabilities on Sheet1 to reflect the re-simulated
# game counts (added more games, sped up simulate game logic, drafted
# optimization logic). Only the numeric probability cells changed; row/column
# labels, styles, and zero-valued cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 0.2185792349726776
$ws.Cells.Item(2, 3).Value = 0.5273224043715847
$ws.Cells.Item(2, 10).Value = 0.01092896174863388
$ws.Cells.Item(2, 16).Value = 0.1475409836065574
$ws.Cells.Item(2, 19).Value = 0.09562841530054644

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01470588235294118
$ws.Cells.Item(3, 3).Value = 0.0392156862745098
$ws.Cells.Item(3, 10).Value = 0.01470588235294118
$ws.Cells.Item(3, 16).Value = 0.7352941176470589
$ws.Cells.Item(3, 19).Value = 0.196078431372549

# Row 4
$ws.Cells.Item(4, 10).Value = 0.05555555555555555
$ws.Cells.Item(4, 16).Value = 0.6296296296296297
$ws.Cells.Item(4, 19).Value = 0.3148148148148148

# Row 6
$ws.Cells.Item(6, 2).Value = 0.06772908366533864
$ws.Cells.Item(6, 4).Value = 0.0199203187250996
$ws.Cells.Item(6, 5).Value = 0.00398406374501992
$ws.Cells.Item(6, 6).Value = 0.05577689243027888
$ws.Cells.Item(6, 10).Value = 0.3147410358565737
$ws.Cells.Item(6, 15).Value = 0.03187250996015936
$ws.Cells.Item(6, 17).Value = 0.1872509960159363
$ws.Cells.Item(6, 18).Value = 0.04382470119521913
$ws.Cells.Item(6, 19).Value = 0.2749003984063745

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1161290322580645
$ws.Cells.Item(7, 4).Value = 0.03225806451612903
$ws.Cells.Item(7, 6).Value = 0.05806451612903226
$ws.Cells.Item(7, 10).Value = 0.1354838709677419
$ws.Cells.Item(7, 15).Value = 0.01290322580645161
$ws.Cells.Item(7, 17).Value = 0.2
$ws.Cells.Item(7, 18).Value = 0.08387096774193549
$ws.Cells.Item(7, 19).Value = 0.3612903225806451

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1587743732590529
$ws.Cells.Item(8, 4).Value = 0.01671309192200557
$ws.Cells.Item(8, 6).Value = 0.08077994428969359
$ws.Cells.Item(8, 10).Value = 0.1030640668523677
$ws.Cells.Item(8, 15).Value = 0.02785515320334262
$ws.Cells.Item(8, 17).Value = 0.2089136490250696
$ws.Cells.Item(8, 18).Value = 0.06963788300835655
$ws.Cells.Item(8, 19).Value = 0.3342618384401114

# Row 9
$ws.Cells.Item(9, 2).Value = 0.08854166666666667
$ws.Cells.Item(9, 4).Value = 0.03125
$ws.Cells.Item(9, 6).Value = 0.0625
$ws.Cells.Item(9, 10).Value = 0.15625
$ws.Cells.Item(9, 15).Value = 0.015625
$ws.Cells.Item(9, 17).Value = 0.2135416666666667
$ws.Cells.Item(9, 18).Value = 0.09895833333333333

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1385245901639344
$ws.Cells.Item(10, 4).Value = 0.02622950819672131
$ws.Cells.Item(10, 5).Value = 0.000819672131147541
$ws.Cells.Item(10, 6).Value = 0.07868852459016394
$ws.Cells.Item(10, 10).Value = 0.1049180327868852
$ws.Cells.Item(10, 15).Value = 0.02786885245901639
$ws.Cells.Item(10, 17).Value = 0.2213114754098361
$ws.Cells.Item(10, 18).Value = 0.07377049180327869
$ws.Cells.Item(10, 19).Value = 0.3278688524590164

# Row 11
$ws.Cells.Item(11, 6).Value = 0.003558718861209964
$ws.Cells.Item(11, 7).Value = 0.1637010676156584
$ws.Cells.Item(11, 10).Value = 0.103202846975089
$ws.Cells.Item(11, 11).Value = 0.2241992882562278
$ws.Cells.Item(11, 12).Value = 0.4911032028469751
$ws.Cells.Item(11, 19).Value = 0.01423487544483986

# Row 12
$ws.Cells.Item(12, 7).Value = 0.6620689655172414
$ws.Cells.Item(12, 10).Value = 0.2413793103448276
$ws.Cells.Item(12, 11).Value = 0.006896551724137931
$ws.Cells.Item(12, 12).Value = 0.02758620689655172
$ws.Cells.Item(12, 19).Value = 0.06206896551724138

# Row 13
$ws.Cells.Item(13, 7).Value = 0.6206896551724138
$ws.Cells.Item(13, 10).Value = 0.3103448275862069
$ws.Cells.Item(13, 19).Value = 0.06896551724137931

# Row 15
$ws.Cells.Item(15, 6).Value = 0.02991452991452992
$ws.Cells.Item(15, 8).Value = 0.1196581196581197
$ws.Cells.Item(15, 9).Value = 0.05128205128205128
$ws.Cells.Item(15, 10).Value = 0.3803418803418803
$ws.Cells.Item(15, 11).Value = 0.04273504273504274
$ws.Cells.Item(15, 13).Value = 0.004273504273504274
$ws.Cells.Item(15, 14).Value = 0.004273504273504274
$ws.Cells.Item(15, 15).Value = 0.07692307692307693
$ws.Cells.Item(15, 19).Value = 0.2905982905982906

# Row 16
$ws.Cells.Item(16, 6).Value = 0.01746724890829694
$ws.Cells.Item(16, 8).Value = 0.1921397379912664
$ws.Cells.Item(16, 9).Value = 0.08296943231441048
$ws.Cells.Item(16, 10).Value = 0.3973799126637554
$ws.Cells.Item(16, 11).Value = 0.1310043668122271
$ws.Cells.Item(16, 13).Value = 0.01310043668122271
$ws.Cells.Item(16, 15).Value = 0.06550218340611354
$ws.Cells.Item(16, 19).Value = 0.1004366812227074

# Row 17
$ws.Cells.Item(17, 6).Value = 0.03463203463203463
$ws.Cells.Item(17, 8).Value = 0.1471861471861472
$ws.Cells.Item(17, 9).Value = 0.119047619047619
$ws.Cells.Item(17, 10).Value = 0.4069264069264069
$ws.Cells.Item(17, 11).Value = 0.09523809523809523
$ws.Cells.Item(17, 13).Value = 0.01731601731601732
$ws.Cells.Item(17, 15).Value = 0.08008658008658008
$ws.Cells.Item(17, 19).Value = 0.09956709956709957

# Row 18
$ws.Cells.Item(18, 6).Value = 0.02531645569620253
$ws.Cells.Item(18, 8).Value = 0.1582278481012658
$ws.Cells.Item(18, 9).Value = 0.06962025316455696
$ws.Cells.Item(18, 10).Value = 0.4430379746835443
$ws.Cells.Item(18, 11).Value = 0.1518987341772152
$ws.Cells.Item(18, 13).Value = 0.006329113924050633
$ws.Cells.Item(18, 15).Value = 0.03164556962025317
$ws.Cells.Item(18, 19).Value = 0.1139240506329114

# Row 19
$ws.Cells.Item(19, 6).Value = 0.03333333333333333
$ws.Cells.Item(19, 8).Value = 0.1851851851851852
$ws.Cells.Item(19, 9).Value = 0.08981481481481482
$ws.Cells.Item(19, 10).Value = 0.3851851851851852
$ws.Cells.Item(19, 11).Value = 0.09722222222222222
$ws.Cells.Item(19, 13).Value = 0.01388888888888889
$ws.Cells.Item(19, 14).Value = 0.000925925925925926
$ws.Cells.Item(19, 15).Value = 0.07314814814814814
$ws.Cells.Item(19, 19).Value = 0.1212962962962963
